$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Remove the "P405D" accessory test-data row (row 10); all rows below
# (Pro32xD, MX2-100, P885D) shift up by one.
$ws.Rows.Item(10).Delete()

# Reflect the post-edit selection (the row that now occupies position 10).
$ws.Range("A10:XFD10").Select()
